$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row (weekly Fruta/Hortaliza price record) needs to be inserted
# right before the existing row 94, pushing rows 94:113 down to 95:114.
$ws.Rows.Item(94).Insert()

# Fill in the newly inserted row 94 with the new record's data. Columns
# A, B, C, E, F, G, H, I, J are identical for every row in this sheet
# (same Mercado / Región / Producto), so write the same literal values
# used throughout the rest of the table to keep it internally consistent.
$ws.Range("A94").Value = 11
$ws.Range("B94").Value = "Vega Monumental Concepción"
$ws.Range("C94").Value = "Bíobío"
$ws.Range("D94").Value = 44736
$ws.Range("E94").Value = 8
$ws.Range("F94").Value = "Fruta"
$ws.Range("G94").Value = 100102
$ws.Range("H94").Value = "Cítricos"
$ws.Range("I94").Value = 100102004
$ws.Range("J94").Value = "Mandarina"
$ws.Range("K94").Value = "Clemenuless"
$ws.Range("L94").Value = "Primera"
$ws.Range("M94").Value = 300
$ws.Range("N94").Value = 6000
$ws.Range("O94").Value = 6500
$ws.Range("P94").Value = 6250
$ws.Range("Q94").Value = "$/caja 18 kilos"
$ws.Range("R94").Value = "Provincia de Limarí"
$ws.Range("S94").Value = 347
$ws.Range("T94").Value = 18
